# "UK Aces, ITA & JAPs soon"
# Adds Italian aces (D41:D45) and Japanese aces (N4:N9) rows, with their
# Wikipedia hyperlinks, and updates the active window's selection/scroll.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Add-AceHyperlink($ws, $cellRef, $name, $url, $tooltip) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $name
    if ($tooltip) {
        $ws.Hyperlinks.Add($rng, $url, "", $tooltip, $url)
    } else {
        $ws.Hyperlinks.Add($rng, $url, "", "", $url)
    }
    # Hyperlinks.Add() stamps the cell text with the "TextToDisplay" value
    # (the URL, matching the existing rows' hyperlink display text) - put
    # the person's name back as the visible cell value afterwards.
    $rng.Value = $name
    $rng.Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# New Italian aces -> D41:D45
# ---------------------------------------------------------------------
Add-AceHyperlink $ws "D41" "Teresio Vittorio Martinoli" "https://en.wikipedia.org/wiki/Teresio_Vittorio_Martinoli" $null
Add-AceHyperlink $ws "D42" "Leonardo Ferrulli"          "https://en.wikipedia.org/wiki/Leonardo_Ferrulli"          $null
Add-AceHyperlink $ws "D43" "Franco Lucchini"            "https://en.wikipedia.org/wiki/Franco_Lucchini"            $null
Add-AceHyperlink $ws "D44" "Franco Bordoni"             "https://en.wikipedia.org/wiki/Franco_Bordoni"             $null
Add-AceHyperlink $ws "D45" "Luigi Gorrini"              "https://en.wikipedia.org/wiki/Luigi_Gorrini"              $null

# ---------------------------------------------------------------------
# New Japanese aces -> N4:N9
# ---------------------------------------------------------------------
Add-AceHyperlink $ws "N4" "Tetsuzo Iwamoto"    "https://en.wikipedia.org/wiki/Tetsuzo_Iwamoto"   $null
Add-AceHyperlink $ws "N5" "Shigeo Fukumoto"    "https://en.wikipedia.org/w/index.php?title=Shigeo_Fukumoto&action=edit&redlink=1" "Shigeo Fukumoto (page does not exist)"
Add-AceHyperlink $ws "N6" "Shoichi Sugita"     "https://en.wikipedia.org/w/index.php?title=Shoichi_Sugita&action=edit&redlink=1"  "Shoichi Sugita (page does not exist)"
Add-AceHyperlink $ws "N7" "Hiromichi Shinohara" "https://en.wikipedia.org/wiki/Hiromichi_Shinohara" $null
Add-AceHyperlink $ws "N8" "Takeo Okumura"       "https://en.wikipedia.org/wiki/Takeo_Okumura"       $null
Add-AceHyperlink $ws "N9" "Satoru Anabuki"      "https://en.wikipedia.org/wiki/Satoru_Anabuki"      $null

# ---------------------------------------------------------------------
# Window / view state: scroll so row 7 is at top, select D34, minimize
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("D34").Select()
$win.WindowState = -4140
